# Weekly data refresh: insert the latest week's price record for
# "Terminal La Palmera de La Serena - Zanahoria" right before the
# existing row 158, pushing all the historical rows (old 158..233)
# down by one (new 159..234).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 158 - this shifts rows 158:233 down to 159:234
# and grows the sheet's used range from A1:R233 to A1:R234.
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A158").Value = 8
$ws.Range("B158").Value = "Terminal La Palmera de La Serena"
$ws.Range("C158").Value = "Coquimbo"
$ws.Range("D158").Value = 44523
$ws.Range("E158").Value = 4
$ws.Range("F158").Value = 100114013
$ws.Range("G158").Value = "Zanahoria"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 600
$ws.Range("K158").Value = 6500
$ws.Range("L158").Value = 7000
$ws.Range("M158").Value = 6750
$ws.Range("N158").Value = "`$/saco 20 kilos"
$ws.Range("O158").Value = "Provincia del Elquí"
$ws.Range("P158").Value = 338
$ws.Range("Q158").Value = 20
$ws.Range("R158").Value = "Hortaliza"
